$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "حازم علاء عبدالعزيز عبدالرحيم"
$ws.Range("B2").Value = "hazemalaa700@gmail.com"
$ws.Range("C2").Value = "https://github.com/Hazemalaa1/Security-Task"

$ws.Range("C2").Select()
